$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card2")

# --- Step 1: add the new event as row 16 ---
# Duplicate row 15 (card/tone-range numbers carried down: card=2, Tones=55,
# Date=1111, Serviced by=ححح) into a fresh row 16 for the new event, while
# row 15 still has its original not-yet-completed blank cells - this way the
# new row mirrors exactly the same (mostly blank) shape row 15 started from.
$ws.Range("A15:O15").Copy()
$ws.Range("A16:O16").PasteSpecial(-4163)   # xlPasteValues - values only, keeps text typing, no style copied

# --- Step 2: complete row 15 ---
# Row 15 was left with several truly-blank cells. Bring it in line with
# every other row on this sheet by filling those blanks with the sheet's
# standard placeholder text "nan" (the values already present - D15=55,
# L15=1111, O15=ححح - are left untouched).
$ws.Range("B15").Value = "nan"
$ws.Range("C15").Value = "nan"
$ws.Range("E15").Value = "nan"
$ws.Range("F15").Value = "nan"
$ws.Range("G15").Value = "nan"
$ws.Range("H15").Value = "nan"
$ws.Range("I15").Value = "nan"
$ws.Range("J15").Value = "nan"
$ws.Range("K15").Value = "nan"
$ws.Range("M15").Value = "nan"
$ws.Range("N15").Value = "nan"
